# Update scenario projections and shiny app
# The scenario table shrinks from 10 rows (Scenario A-H plus header/counterfactual
# rows) down to 6 rows (header, counterfactual row, Scenario A-D), and the
# wording of several labels is refreshed with explicit percentages.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-unused Scenario E-H rows (old rows 7-10) first so the
# remaining writes land on the final A1:D6 grid.
$ws.Rows("7:10").Delete()

# Row 1 header wording is unchanged.
$ws.Range("A1").Value = "Scenario"
$ws.Range("B1").Value = "Monoclonal Antibodies "
$ws.Range("C1").Value = "Maternal Vaccination"
$ws.Range("D1").Value = "Senior Vaccination"

# Write the brand-new label text in the order it is first introduced through
# the refreshed table, so freshly interned shared strings line up with how
# Excel appends them to the shared string table.
$ws.Range("A2").Value = "ScenarioE (Counterfactual)"
$ws.Range("D2").Value = "25% (from 2023-34 season)"
$ws.Range("D3").Value = "Optimistic (40%)"
$ws.Range("D4").Value = "Pessimistic (30%)"
$ws.Range("B3").Value = "Optimistic (75%)"
$ws.Range("B5").Value = "Pessimistic (25%)"
$ws.Range("C3").Value = "Optimistic (50%)"

# Remaining cells re-use the labels already introduced above.
$ws.Range("B2").Value = "None"
$ws.Range("C2").Value = "None"

$ws.Range("A3").Value = "Scenario A"
$ws.Range("C4").Value = "Optimistic (50%)"

$ws.Range("A4").Value = "Scenario B"
$ws.Range("B4").Value = "Optimistic (75%)"

$ws.Range("A5").Value = "Scenario C"
$ws.Range("C5").Value = "Pessimistic (30%)"
$ws.Range("D5").Value = "Optimistic (40%)"

$ws.Range("A6").Value = "Scenario D"
$ws.Range("B6").Value = "Pessimistic (25%)"
$ws.Range("C6").Value = "Pessimistic (30%)"
$ws.Range("D6").Value = "Pessimistic (30%)"

# Selection moved to F8 in the saved view
$ws.Range("F8").Select()
